$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the crypto price/volume refresh.
# Values that would otherwise be auto-parsed as numbers by Excel
# are written with a leading quote (forcing text entry, matching
# the original cells' text/General storage) and then restored to
# the default "Normal" style so no stray number formatting remains.

$ws.Range('D2').Value = '65.333.65'
$ws.Range('E2').Value = '  -3.68%  '
$ws.Range('D3').Value = '3.469.36'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'562.19"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').Value = "'175.14"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.63%  '
$ws.Range('D7').Value = "'0.625"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.70%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = "'0.624"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('D10').Value = "'0.152"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('D11').Value = "'53.41"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.39%  '
$ws.Range('D12').Value = "'0.0000268"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.26%  '
$ws.Range('D13').Value = "'9.08"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('D14').Value = '4.033.27'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = '3.484.03'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').Value = "'18.11"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = "'11.99"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = '65.382.40'
$ws.Range('E19').Value = '  -3.60%  '
$ws.Range('D20').Value = "'0.987"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').Value = "'409.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = "'4.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.62%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = "'4.27"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'84.63"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').Value = "'12.72"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.28%  '
$ws.Range('D26').Value = "'10.73"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.60%  '
$ws.Range('D27').Value = "'2.80"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.13%  '
$ws.Range('D28').Value = "'8.86"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('D29').Value = "'29.92"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').Value = "'620.44"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -11.53%  '
$ws.Range('D31').Value = "'6.28"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.40%  '
$ws.Range('D32').Value = "'11.50"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').Value = "'0.108"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.48%  '
$ws.Range('D34').Value = "'58.69"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('D35').Value = "'0.146"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.77%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '0.0₃0782'
$ws.Range('E37').Value = '  -6.22%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '3.357.68'
$ws.Range('E38').Value = '  +10.65%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = "'36.58"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.95%  '
$ws.Range('D40').Value = "'0.374"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.12%  '
$ws.Range('D41').Value = "'3.40"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').Value = "'2.85"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.61%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = "'3.24"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.03%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = "'0.0410"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.07%  '
$ws.Range('D46').Value = "'2.47"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.68%  '
$ws.Range('E47').Value = '  -1.88%  '
$ws.Range('D48').Value = "'0.131"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('D49').Value = "'136.89"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('D50').Value = "'8.32"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.23%  '
$ws.Range('D51').Value = "'2.79"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.54%  '
